$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 0.00006240767534437808
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1133.908754285003
